$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Object_Mapping")

# Fix spelling: Destilation -> Distillation
$ws.Range("A5").Value = "Distillation_Tower"
$ws.Range("B5").Value = "Distillation_tower"

# Abbreviate pipeline/power_line -> PL
$ws.Range("A8").Value = "PL_Wholesale_Kasso"
$ws.Range("A9").Value = "PL_storage_hydrogen"
$ws.Range("A10").Value = "PL_storage_e-methanol"
$ws.Range("A11").Value = "PL_District_Heating"

# Drop "_Kasso" suffix from storage object names
$ws.Range("A12").Value = "Hydrogen_storage"
$ws.Range("A13").Value = "E-Methanol_storage"
